$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove all existing hyperlinks so the set below can be rebuilt cleanly ---
$ws.Cells.Hyperlinks.Delete()

# --- 2. Update the data area text ---
# Header row stays the same text
$ws.Range("A1").Value = "Email Id"
$ws.Range("B1").Value = "Password"

# New set of valid-login test credentials
$ws.Range("A2").Value = "rafselenium1@gmail.com"
$ws.Range("B2").Value = "Selenium@123"

$ws.Range("A3").Value = "rafselenium2@gmail.com"
$ws.Range("B3").Value = "Selenium@123"

$ws.Range("A4").Value = "rafselenium3@yahoo.com"
$ws.Range("B4").Value = "Selenium@123"

$ws.Range("A5").Value = "rafselenium4@yahoo.com"
$ws.Range("B5").Value = "Selenium@123"

# Row 6 no longer holds any text, just blank styled cells
$ws.Range("A6:B6").ClearContents()

# --- 3. Re-create the hyperlinks (same mailto targets as before; the one that used
#        to sit on B6 now sits on A4 since A4 carries a real address now) ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:rafiasultana12345@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Selenium@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:rafiasultana122@yahoo.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Selenium@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:sharmin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:rashidmohammed@yahoo.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:mohammed@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:rashida@123") | Out-Null

# --- 4. Style the blank cells (row 6 and the 5 new blank rows) with the same
#        hyperlink-colored, border-less look ---
$ws.Range("A6:B6").Style = "Hyperlink"
$ws.Range("A9:B13").Style = "Hyperlink"

# --- 5. Selection / active cell moves to A9 ---
$ws.Range("A9").Select()
